# Applies the NATMI TPM data refresh for the Lgi2-Adam22 sheet (rows 2-19).
# For each row the Target cluster label (col D) and the numeric result columns
# (E,F,G,H,I,J,M,N,O,P,Q,R,S,T) are updated to the newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.3777555
$ws.Range("H2").Value = 0.755511
$ws.Range("I2").Value = 0.04674878051708328
$ws.Range("J2").Value = 0.03234814298672928
$ws.Range("M2").Value = 5.447678
$ws.Range("N2").Value = 10.895356
$ws.Range("O2").Value = 0.4286498436662743
$ws.Range("P2").Value = 0.4047900009176674
$ws.Range("Q2").Value = 2.057890326729
$ws.Range("R2").Value = 8.231561306916
$ws.Range("S2").Value = 0.02003885746023672
$ws.Range("T2").Value = 0.01309420482928298

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.3777555
$ws.Range("H3").Value = 0.755511
$ws.Range("I3").Value = 0.04674878051708328
$ws.Range("J3").Value = 0.03234814298672928
$ws.Range("M3").Value = 0.7327576666666668
$ws.Range("O3").Value = 0.05765694287766837
$ws.Range("P3").Value = 0.08167139556406268
$ws.Range("Q3").Value = 0.2768032387505001
$ws.Range("R3").Value = 1.660819432503
$ws.Range("S3").Value = 0.002695391767874127
$ws.Range("T3").Value = 0.002641917981632028

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.3777555
$ws.Range("H4").Value = 0.755511
$ws.Range("I4").Value = 0.04674878051708328
$ws.Range("J4").Value = 0.03234814298672928
$ws.Range("M4").Value = 0.302684
$ws.Range("N4").Value = 0.9080520000000001
$ws.Range("O4").Value = 0.02381665165971311
$ws.Range("P4").Value = 0.03373642586009028
$ws.Range("Q4").Value = 0.114340545762
$ws.Range("R4").Value = 0.6860432745720001
$ws.Range("S4").Value = 0.001113399421091755
$ws.Range("T4").Value = 0.001091310727583392

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 0.3777555
$ws.Range("H5").Value = 0.755511
$ws.Range("I5").Value = 0.04674878051708328
$ws.Range("J5").Value = 0.03234814298672928
$ws.Range("M5").Value = 5.763022
$ws.Range("N5").Value = 11.526044
$ws.Range("O5").Value = 0.4534626457997884
$ws.Range("P5").Value = 0.4282216534583244
$ws.Range("Q5").Value = 2.177013257121
$ws.Range("R5").Value = 8.708053028484
$ws.Range("S5").Value = 0.02119882570119019
$ws.Range("T5").Value = 0.01385217527608351

# Row 6
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 0.3777555
$ws.Range("H6").Value = 0.755511
$ws.Range("I6").Value = 0.04674878051708328
$ws.Range("J6").Value = 0.03234814298672928
$ws.Range("M6").Value = 0.2632226666666667
$ws.Range("N6").Value = 0.789668
$ws.Range("O6").Value = 0.02071164171525676
$ws.Range("P6").Value = 0.02933816118029118
$ws.Range("Q6").Value = 0.099433810058
$ws.Range("R6").Value = 0.596602860348
$ws.Range("S6").Value = 0.0009682439926950046
$ws.Range("T6").Value = 0.0009490350328277696

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 0.3777555
$ws.Range("H7").Value = 0.755511
$ws.Range("I7").Value = 0.04674878051708328
$ws.Range("J7").Value = 0.03234814298672928
$ws.Range("M7").Value = 0.199559
$ws.Range("N7").Value = 0.598677
$ws.Range("O7").Value = 0.01570227428129894
$ws.Range("P7").Value = 0.02224236301956415
$ws.Range("Q7").Value = 0.07538450982450001
$ws.Range("R7").Value = 0.452307058947
$ws.Range("S7").Value = 0.0007340621739954858
$ws.Range("T7").Value = 0.0007194991393196008

# Row 8
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.194537333333333
$ws.Range("H8").Value = 21.583612
$ws.Range("I8").Value = 0.8903532753804024
$ws.Range("J8").Value = 0.9241291882528327
$ws.Range("M8").Value = 5.447678
$ws.Range("N8").Value = 10.895356
$ws.Range("O8").Value = 0.4286498436662743
$ws.Range("P8").Value = 0.4047900009176674
$ws.Range("Q8").Value = 39.19352275097867
$ws.Range("R8").Value = 235.161136505872
$ws.Range("S8").Value = 0.3816497922995648
$ws.Range("T8").Value = 0.3740782549609074

# Row 9
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.194537333333333
$ws.Range("H9").Value = 21.583612
$ws.Range("I9").Value = 0.8903532753804024
$ws.Range("J9").Value = 0.9241291882528327
$ws.Range("M9").Value = 0.7327576666666668
$ws.Range("O9").Value = 0.05765694287766837
$ws.Range("P9").Value = 0.08167139556406268
$ws.Range("Q9").Value = 5.271852389119556
$ws.Range("R9").Value = 47.44667150207601
$ws.Range("S9").Value = 0.0513350479395528
$ws.Range("T9").Value = 0.07547492048609325

# Row 10
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.194537333333333
$ws.Range("H10").Value = 21.583612
$ws.Range("I10").Value = 0.8903532753804024
$ws.Range("J10").Value = 0.9241291882528327
$ws.Range("M10").Value = 0.302684
$ws.Range("N10").Value = 0.9080520000000001
$ws.Range("O10").Value = 0.02381665165971311
$ws.Range("P10").Value = 0.03373642586009028
$ws.Range("Q10").Value = 2.177671338202666
$ws.Range("R10").Value = 19.599042043824
$ws.Range("S10").Value = 0.02120523381381966
$ws.Range("T10").Value = 0.0311768158446371

# Row 11
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.194537333333333
$ws.Range("H11").Value = 21.583612
$ws.Range("I11").Value = 0.8903532753804024
$ws.Range("J11").Value = 0.9241291882528327
$ws.Range("M11").Value = 5.763022
$ws.Range("N11").Value = 11.526044
$ws.Range("O11").Value = 0.4534626457997884
$ws.Range("P11").Value = 0.4282216534583244
$ws.Range("Q11").Value = 41.46227693182134
$ws.Range("R11").Value = 248.773661590928
$ws.Range("S11").Value = 0.4037419519505049
$ws.Range("T11").Value = 0.3957321290027272

# Row 12
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 7.194537333333333
$ws.Range("H12").Value = 21.583612
$ws.Range("I12").Value = 0.8903532753804024
$ws.Range("J12").Value = 0.9241291882528327
$ws.Range("M12").Value = 0.2632226666666667
$ws.Range("N12").Value = 0.789668
$ws.Range("O12").Value = 0.02071164171525676
$ws.Range("P12").Value = 0.02933816118029118
$ws.Range("Q12").Value = 1.893765302312889
$ws.Range("R12").Value = 17.043887720816
$ws.Range("S12").Value = 0.01844067803968423
$ws.Range("T12").Value = 0.02711225107637326

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 7.194537333333333
$ws.Range("H13").Value = 21.583612
$ws.Range("I13").Value = 0.8903532753804024
$ws.Range("J13").Value = 0.9241291882528327
$ws.Range("M13").Value = 0.199559
$ws.Range("N13").Value = 0.598677
$ws.Range("O13").Value = 0.01570227428129894
$ws.Range("P13").Value = 0.02224236301956415
$ws.Range("Q13").Value = 1.435734675702667
$ws.Range("R13").Value = 12.921612081324
$ws.Range("S13").Value = 0.01398057133727597
$ws.Range("T13").Value = 0.02055481688209464

# Row 14
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5082495
$ws.Range("H14").Value = 1.016499
$ws.Range("I14").Value = 0.06289794410251424
$ws.Range("J14").Value = 0.04352266876043807
$ws.Range("M14").Value = 5.447678
$ws.Range("N14").Value = 10.895356
$ws.Range("O14").Value = 0.4286498436662743
$ws.Range("P14").Value = 0.4047900009176674
$ws.Range("Q14").Value = 2.768779619661
$ws.Range("R14").Value = 11.075118478644
$ws.Range("S14").Value = 0.02696119390647279
$ws.Range("T14").Value = 0.01761754112747706

# Row 15
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5082495
$ws.Range("H15").Value = 1.016499
$ws.Range("I15").Value = 0.06289794410251424
$ws.Range("J15").Value = 0.04352266876043807
$ws.Range("M15").Value = 0.7327576666666668
$ws.Range("O15").Value = 0.05765694287766837
$ws.Range("P15").Value = 0.08167139556406268
$ws.Range("Q15").Value = 0.3724237177045001
$ws.Range("R15").Value = 2.234542306227
$ws.Range("S15").Value = 0.003626503170241442
$ws.Range("T15").Value = 0.003554557096337411

# Row 16
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5082495
$ws.Range("H16").Value = 1.016499
$ws.Range("I16").Value = 0.06289794410251424
$ws.Range("J16").Value = 0.04352266876043807
$ws.Range("M16").Value = 0.302684
$ws.Range("N16").Value = 0.9080520000000001
$ws.Range("O16").Value = 0.02381665165971311
$ws.Range("P16").Value = 0.03373642586009028
$ws.Range("Q16").Value = 0.153838991658
$ws.Range("R16").Value = 0.9230339499480001
$ws.Range("S16").Value = 0.001498018424801688
$ws.Range("T16").Value = 0.001468299287869786

# Row 17
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5082495
$ws.Range("H17").Value = 1.016499
$ws.Range("I17").Value = 0.06289794410251424
$ws.Range("J17").Value = 0.04352266876043807
$ws.Range("M17").Value = 5.763022
$ws.Range("N17").Value = 11.526044
$ws.Range("O17").Value = 0.4534626457997884
$ws.Range("P17").Value = 0.4282216534583244
$ws.Range("Q17").Value = 2.929053049989
$ws.Range("R17").Value = 11.716212199956
$ws.Range("S17").Value = 0.02852186814809331
$ws.Range("T17").Value = 0.01863734917951376

# Row 18
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.5082495
$ws.Range("H18").Value = 1.016499
$ws.Range("I18").Value = 0.06289794410251424
$ws.Range("J18").Value = 0.04352266876043807
$ws.Range("M18").Value = 0.2632226666666667
$ws.Range("N18").Value = 0.789668
$ws.Range("O18").Value = 0.02071164171525676
$ws.Range("P18").Value = 0.02933816118029118
$ws.Range("Q18").Value = 0.133782788722
$ws.Range("R18").Value = 0.8026967323320001
$ws.Range("S18").Value = 0.001302719682877522
$ws.Range("T18").Value = 0.001276875071090156

# Row 19
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.5082495
$ws.Range("H19").Value = 1.016499
$ws.Range("I19").Value = 0.06289794410251424
$ws.Range("J19").Value = 0.04352266876043807
$ws.Range("M19").Value = 0.199559
$ws.Range("N19").Value = 0.598677
$ws.Range("O19").Value = 0.01570227428129894
$ws.Range("P19").Value = 0.02224236301956415
$ws.Range("Q19").Value = 0.1014257619705
$ws.Range("R19").Value = 0.608554571823
$ws.Range("S19").Value = 0.0009876407700274879
$ws.Range("T19").Value = 0.0009680469981499077
